$wb = $excel.ActiveWorkbook

$changes = @{
    "F2"  = 365
    "F4"  = 10799
    "F6"  = 976
    "F7"  = 164
    "F8"  = 1336
    "F9"  = 8289
    "F10" = 38
    "F11" = 467
    "F12" = 422
    "F15" = 3301
    "F16" = 40
    "F17" = 327
    "F18" = 26
    "F19" = 782
    "F21" = 1073
    "F23" = 109
    "F24" = 1772
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $changes.Keys) {
        $ws.Range($addr).Value = $changes[$addr]
    }
}
